# The edit swaps the full data of row 16 and row 17 (two species
# observation records trade places in the sheet), including the stray
# empty "Kön" cell in column L that belongs with the Grönpyrola / Pyrola
# chlorantha record.
#
# Columns that actually carry different values between the two records:
#   A  Id
#   B  Taxonsorteringsordning
#   E  TaxonId
#   F  Artnamn
#   G  Vetenskapligt namn
#   H  Auktor
#   J  Enhet
#   Q  Ost
#   R  Nord
# plus the empty column-L cell, which moves from row 16 to row 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$swapCols = @("A", "B", "E", "F", "G", "H", "J", "Q", "R")

# Stash row 16's current values for the swapping columns in a scratch
# area far below the used range, then pull row 17's values up into row
# 16, and finally drop the stashed values into row 17. Using
# Range.Copy (rather than .Value assignment) preserves each cell's
# original type/formatting instead of collapsing everything to text.
foreach ($col in $swapCols) {
    $ws.Range($col + "16").Copy($ws.Range($col + "9016"))
}
foreach ($col in $swapCols) {
    $ws.Range($col + "17").Copy($ws.Range($col + "16"))
}
foreach ($col in $swapCols) {
    $ws.Range($col + "9016").Copy($ws.Range($col + "17"))
}

# Clear the scratch row so it doesn't linger in the saved sheet.
$ws.Range("A9016:R9016").ClearContents()

# Column L (an empty "Kön" cell) belongs to the record that is now in
# row 17 (previously row 16), so move it along with the rest of the
# data: create it in row 17 and remove it from row 16.
$ws.Range("L16").Copy($ws.Range("L17"))
$ws.Range("L16").ClearContents()
